$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 1313.3846  # ALC!H2 was 1402.4166
$ws.Cells.Item(2, 9).Value = 1747.2222  # ALC!I2 was 1935
$ws.Cells.Item(2, 11).Value = 1747.2222  # ALC!K2 was 1935
$ws.Cells.Item(2, 13).Value = -1634.2222  # ALC!M2 was -1822
$ws.Cells.Item(6, 8).Value = 6243  # ALC!H6 was 12385.714
$ws.Cells.Item(6, 9).Value = 7280.5835  # ALC!I6 was 12385.714
$ws.Cells.Item(6, 10).Value = 17.5  # ALC!J6 was 0
$ws.Cells.Item(6, 11).Value = 21841.7505  # ALC!K6 was 37157.142
$ws.Cells.Item(6, 12).Value = 52.5  # ALC!L6 was 0
$ws.Cells.Item(6, 13).Value = -21729.7505  # ALC!M6 was -37045.142
$ws.Cells.Item(6, 14).Value = -276.5  # ALC!N6 (new cell)
$ws.Cells.Item(80, 8).Value = 541.0833  # ALC!H80 was 549.4167
$ws.Cells.Item(80, 9).Value = 274.75  # ALC!I80 was 269.8
$ws.Cells.Item(80, 10).Value = 674.25  # ALC!J80 was 749.1429000000001
$ws.Cells.Item(80, 11).Value = 824.25  # ALC!K80 was 809.4000000000001
$ws.Cells.Item(80, 12).Value = 2022.75  # ALC!L80 was 2247.4287
$ws.Cells.Item(80, 13).Value = 173.75  # ALC!M80 was 188.5999999999999
$ws.Cells.Item(80, 14).Value = -4018.75  # ALC!N80 was -4243.4287
$ws.Cells.Item(83, 8).Value = 541.0833  # ALC!H83 was 549.4167
$ws.Cells.Item(83, 9).Value = 274.75  # ALC!I83 was 269.8
$ws.Cells.Item(83, 10).Value = 674.25  # ALC!J83 was 749.1429000000001
$ws.Cells.Item(83, 11).Value = 2472.75  # ALC!K83 was 2428.2
$ws.Cells.Item(83, 12).Value = 6068.25  # ALC!L83 was 6742.2861
$ws.Cells.Item(83, 13).Value = 2519.25  # ALC!M83 was 2563.8
$ws.Cells.Item(83, 14).Value = -16052.25  # ALC!N83 was -16726.2861
$ws.Cells.Item(92, 8).Value = 205.05882  # ALC!H92 was 215.93333
$ws.Cells.Item(92, 9).Value = 206.3125  # ALC!I92 was 218.14285
$ws.Cells.Item(92, 11).Value = 206.3125  # ALC!K92 was 218.14285
$ws.Cells.Item(92, 13).Value = 1041.6875  # ALC!M92 was 1029.85715
$ws.Cells.Item(99, 8).Value = 517.7368  # ALC!H99 was 490.42105
$ws.Cells.Item(99, 9).Value = 379.92307  # ALC!I99 was 358.5
$ws.Cells.Item(99, 10).Value = 816.3333  # ALC!J99 was 859.8
$ws.Cells.Item(99, 11).Value = 1139.76921  # ALC!K99 was 1075.5
$ws.Cells.Item(99, 12).Value = 2448.9999  # ALC!L99 was 2579.4
$ws.Cells.Item(99, 13).Value = 358.2307900000001  # ALC!M99 was 422.5
$ws.Cells.Item(99, 14).Value = -5444.9999  # ALC!N99 was -5575.4
$ws.Cells.Item(101, 8).Value = 1628.4706  # ALC!H101 was 1635.8235
$ws.Cells.Item(101, 9).Value = 1946.6  # ALC!I101 was 1959.1
$ws.Cells.Item(101, 11).Value = 5839.799999999999  # ALC!K101 was 5877.299999999999
$ws.Cells.Item(101, 13).Value = -4217.799999999999  # ALC!M101 was -4255.299999999999
$ws.Cells.Item(112, 8).Value = 1900.1  # ALC!H112 was 1897.5264
$ws.Cells.Item(112, 10).Value = 2089.111  # ALC!J112 was 2097.353
$ws.Cells.Item(112, 12).Value = 6267.333  # ALC!L112 was 6292.059
$ws.Cells.Item(112, 14).Value = -8483.332999999999  # ALC!N112 was -8508.059000000001
$ws.Cells.Item(132, 8).Value = 2490.6155  # ALC!H132 was 2347.7144
$ws.Cells.Item(132, 9).Value = 2490.6155  # ALC!I132 was 2347.7144
$ws.Cells.Item(132, 11).Value = 7471.8465  # ALC!K132 was 7043.1432
$ws.Cells.Item(132, 13).Value = -4941.8465  # ALC!M132 was -4513.1432
$ws.Cells.Item(137, 8).Value = 2832.0571  # ALC!H137 was 2857.7715
$ws.Cells.Item(137, 9).Value = 3243.7917  # ALC!I137 was 3281.2917
$ws.Cells.Item(137, 11).Value = 9731.375100000001  # ALC!K137 was 9843.875100000001
$ws.Cells.Item(137, 13).Value = -7181.375100000001  # ALC!M137 was -7293.875100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 278.125  # ARM!H5 was 295.46667
$ws.Cells.Item(5, 9).Value = 124.44444  # ARM!I5 was 133.55556
$ws.Cells.Item(5, 10).Value = 475.7143  # ARM!J5 was 538.3333
$ws.Cells.Item(5, 11).Value = 124.44444  # ARM!K5 was 133.55556
$ws.Cells.Item(5, 12).Value = 475.7143  # ARM!L5 was 538.3333
$ws.Cells.Item(5, 13).Value = -12.44444  # ARM!M5 was -21.55556000000001
$ws.Cells.Item(5, 14).Value = -699.7143  # ARM!N5 was -762.3333
$ws.Cells.Item(32, 8).Value = 21756412  # ARM!H32 was 15163785
$ws.Cells.Item(32, 9).Value = 29426344  # ARM!I32 was 18527972
$ws.Cells.Item(32, 11).Value = 29426344  # ARM!K32 was 18527972
$ws.Cells.Item(32, 13).Value = -29426057  # ARM!M32 was -18527685
$ws.Cells.Item(74, 8).Value = 11819118  # ARM!H74 was 11819095
$ws.Cells.Item(74, 9).Value = 14706745  # ARM!I74 was 14706715
$ws.Cells.Item(74, 11).Value = 14706745  # ARM!K74 was 14706715
$ws.Cells.Item(74, 13).Value = -14705871  # ARM!M74 was -14705841
$ws.Cells.Item(76, 8).Value = 23644  # ARM!H76 was 40287.332
$ws.Cells.Item(76, 9).Value = 7000  # ARM!I76 was 0
$ws.Cells.Item(76, 10).Value = 40288  # ARM!J76 was 40287.332
$ws.Cells.Item(76, 11).Value = 7000  # ARM!K76 was 0
$ws.Cells.Item(76, 12).Value = 40288  # ARM!L76 was 40287.332
$ws.Cells.Item(76, 13).Value = -6662  # ARM!M76 (new cell)
$ws.Cells.Item(76, 14).Value = -40964  # ARM!N76 was -40963.332
$ws.Cells.Item(77, 8).Value = 11819118  # ARM!H77 was 11819095
$ws.Cells.Item(77, 9).Value = 14706745  # ARM!I77 was 14706715
$ws.Cells.Item(77, 11).Value = 73533725  # ARM!K77 was 73533575
$ws.Cells.Item(77, 13).Value = -73529357  # ARM!M77 was -73529207
$ws.Cells.Item(79, 8).Value = 23644  # ARM!H79 was 40287.332
$ws.Cells.Item(79, 9).Value = 7000  # ARM!I79 was 0
$ws.Cells.Item(79, 10).Value = 40288  # ARM!J79 was 40287.332
$ws.Cells.Item(79, 11).Value = 7000  # ARM!K79 was 0
$ws.Cells.Item(79, 12).Value = 40288  # ARM!L79 was 40287.332
$ws.Cells.Item(79, 13).Value = -5830  # ARM!M79 (new cell)
$ws.Cells.Item(79, 14).Value = -42628  # ARM!N79 was -42627.332
$ws.Cells.Item(97, 8).Value = 1386.174  # ARM!H97 was 1582.6522
$ws.Cells.Item(97, 9).Value = 861  # ARM!I97 was 1183.7858
$ws.Cells.Item(97, 11).Value = 861  # ARM!K97 was 1183.7858
$ws.Cells.Item(97, 13).Value = -365  # ARM!M97 was -687.7858000000001
$ws.Cells.Item(102, 8).Value = 13413.154  # ARM!H102 was 14468.417
$ws.Cells.Item(102, 9).Value = 13413.154  # ARM!I102 was 14468.417
$ws.Cells.Item(102, 11).Value = 13413.154  # ARM!K102 was 14468.417
$ws.Cells.Item(102, 13).Value = -11791.154  # ARM!M102 was -12846.417
$ws.Cells.Item(103, 8).Value = 74989  # ARM!H103 was 51598
$ws.Cells.Item(103, 10).Value = 74989  # ARM!J103 was 51598
$ws.Cells.Item(103, 12).Value = 74989  # ARM!L103 was 51598
$ws.Cells.Item(103, 14).Value = -77333  # ARM!N103 was -53942

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 278.125  # BSM!H4 was 295.46667
$ws.Cells.Item(4, 9).Value = 124.44444  # BSM!I4 was 133.55556
$ws.Cells.Item(4, 10).Value = 475.7143  # BSM!J4 was 538.3333
$ws.Cells.Item(4, 11).Value = 124.44444  # BSM!K4 was 133.55556
$ws.Cells.Item(4, 12).Value = 475.7143  # BSM!L4 was 538.3333
$ws.Cells.Item(4, 13).Value = -9.44444  # BSM!M4 was -18.55556000000001
$ws.Cells.Item(4, 14).Value = -705.7143  # BSM!N4 was -768.3333
$ws.Cells.Item(108, 8).Value = 104936.336  # BSM!H108 was 104943
$ws.Cells.Item(108, 10).Value = 104936.336  # BSM!J108 was 104943
$ws.Cells.Item(108, 12).Value = 104936.336  # BSM!L108 was 104943
$ws.Cells.Item(108, 14).Value = -112616.336  # BSM!N108 was -112623

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 561682.7  # CRP!H31 was 604732.7
$ws.Cells.Item(31, 9).Value = 9041.261  # CRP!I31 was 9361.362999999999
$ws.Cells.Item(31, 10).Value = 1230669.8  # CRP!J31 was 1375213.2
$ws.Cells.Item(31, 11).Value = 9041.261  # CRP!K31 was 9361.362999999999
$ws.Cells.Item(31, 12).Value = 1230669.8  # CRP!L31 was 1375213.2
$ws.Cells.Item(31, 13).Value = -8746.261  # CRP!M31 was -9066.362999999999
$ws.Cells.Item(31, 14).Value = -1231259.8  # CRP!N31 was -1375803.2
$ws.Cells.Item(34, 8).Value = 561682.7  # CRP!H34 was 604732.7
$ws.Cells.Item(34, 9).Value = 9041.261  # CRP!I34 was 9361.362999999999
$ws.Cells.Item(34, 10).Value = 1230669.8  # CRP!J34 was 1375213.2
$ws.Cells.Item(34, 11).Value = 9041.261  # CRP!K34 was 9361.362999999999
$ws.Cells.Item(34, 12).Value = 1230669.8  # CRP!L34 was 1375213.2
$ws.Cells.Item(34, 13).Value = -8839.261  # CRP!M34 was -9159.362999999999
$ws.Cells.Item(34, 14).Value = -1231073.8  # CRP!N34 was -1375617.2
$ws.Cells.Item(62, 8).Value = 632774.5600000001  # CRP!H62 was 722203.4399999999
$ws.Cells.Item(62, 9).Value = 916782.5600000001  # CRP!I62 was 916869
$ws.Cells.Item(62, 10).Value = 7957  # CRP!J62 was 8429.666999999999
$ws.Cells.Item(62, 11).Value = 916782.5600000001  # CRP!K62 was 916869
$ws.Cells.Item(62, 12).Value = 7957  # CRP!L62 was 8429.666999999999
$ws.Cells.Item(62, 13).Value = -916158.5600000001  # CRP!M62 was -916245
$ws.Cells.Item(62, 14).Value = -9205  # CRP!N62 was -9677.666999999999
$ws.Cells.Item(65, 8).Value = 632774.5600000001  # CRP!H65 was 722203.4399999999
$ws.Cells.Item(65, 9).Value = 916782.5600000001  # CRP!I65 was 916869
$ws.Cells.Item(65, 10).Value = 7957  # CRP!J65 was 8429.666999999999
$ws.Cells.Item(65, 11).Value = 4583912.800000001  # CRP!K65 was 4584345
$ws.Cells.Item(65, 12).Value = 39785  # CRP!L65 was 42148.335
$ws.Cells.Item(65, 13).Value = -4580792.800000001  # CRP!M65 was -4581225
$ws.Cells.Item(65, 14).Value = -46025  # CRP!N65 was -48388.335

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 139.05263  # CUL!H2 was 129.2
$ws.Cells.Item(2, 10).Value = 173.04546  # CUL!J2 was 276.75
$ws.Cells.Item(2, 12).Value = 1038.27276  # CUL!L2 was 1660.5
$ws.Cells.Item(2, 14).Value = -1264.27276  # CUL!N2 was -1886.5
$ws.Cells.Item(11, 8).Value = 130.5  # CUL!H11 was 861.6667
$ws.Cells.Item(11, 9).Value = 130.5  # CUL!I11 was 861.6667
$ws.Cells.Item(11, 11).Value = 391.5  # CUL!K11 was 2585.0001
$ws.Cells.Item(11, 13).Value = -251.5  # CUL!M11 was -2445.0001
$ws.Cells.Item(87, 8).Value = 3299.5715  # CUL!H87 was 2832.4443
$ws.Cells.Item(87, 9).Value = 3299.5715  # CUL!I87 was 2832.4443
$ws.Cells.Item(87, 11).Value = 9898.7145  # CUL!K87 was 8497.332900000001
$ws.Cells.Item(87, 13).Value = -8650.7145  # CUL!M87 was -7249.332900000001
$ws.Cells.Item(90, 8).Value = 3299.5715  # CUL!H90 was 2832.4443
$ws.Cells.Item(90, 9).Value = 3299.5715  # CUL!I90 was 2832.4443
$ws.Cells.Item(90, 11).Value = 29696.1435  # CUL!K90 was 25491.9987
$ws.Cells.Item(90, 13).Value = -23456.1435  # CUL!M90 was -19251.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 0  # GSM!H18 was 59990
$ws.Cells.Item(18, 9).Value = 0  # GSM!I18 was 59990
$ws.Cells.Item(18, 11).Value = 0  # GSM!K18 was 59990
$ws.Cells.Item(18, 13).ClearContents()  # GSM!M18
$ws.Cells.Item(134, 8).Value = 80000  # GSM!H134 was 77499.75
$ws.Cells.Item(134, 10).Value = 80000  # GSM!J134 was 77499.75
$ws.Cells.Item(134, 12).Value = 240000  # GSM!L134 was 232499.25
$ws.Cells.Item(134, 14).Value = -245070  # GSM!N134 was -237569.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 0  # LTW!H68 was 4000
$ws.Cells.Item(68, 9).Value = 0  # LTW!I68 was 4000
$ws.Cells.Item(68, 11).Value = 0  # LTW!K68 was 4000
$ws.Cells.Item(68, 13).ClearContents()  # LTW!M68
$ws.Cells.Item(71, 8).Value = 0  # LTW!H71 was 4000
$ws.Cells.Item(71, 9).Value = 0  # LTW!I71 was 4000
$ws.Cells.Item(71, 11).Value = 0  # LTW!K71 was 20000
$ws.Cells.Item(71, 13).ClearContents()  # LTW!M71

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 12561.875  # WVR!H39 was 13949.4
$ws.Cells.Item(39, 9).Value = 10000  # WVR!I39 was 9875
$ws.Cells.Item(39, 10).Value = 30495  # WVR!J39 was 30247
$ws.Cells.Item(39, 11).Value = 10000  # WVR!K39 was 9875
$ws.Cells.Item(39, 12).Value = 30495  # WVR!L39 was 30247
$ws.Cells.Item(39, 13).Value = -9587  # WVR!M39 was -9462
$ws.Cells.Item(39, 14).Value = -31321  # WVR!N39 was -31073
$ws.Cells.Item(43, 8).Value = 92009  # WVR!H43 was 97348
$ws.Cells.Item(43, 9).Value = 96027  # WVR!I43 was 96022
$ws.Cells.Item(43, 10).Value = 90000  # WVR!J43 was 100000
$ws.Cells.Item(43, 11).Value = 96027  # WVR!K43 was 96022
$ws.Cells.Item(43, 12).Value = 90000  # WVR!L43 was 100000
$ws.Cells.Item(43, 13).Value = -95878  # WVR!M43 was -95873
$ws.Cells.Item(43, 14).Value = -90298  # WVR!N43 was -100298
$ws.Cells.Item(49, 8).Value = 33494.5  # WVR!H49 was 33495
$ws.Cells.Item(49, 9).Value = 33494  # WVR!I49 was 0
$ws.Cells.Item(49, 11).Value = 33494  # WVR!K49 was 0
$ws.Cells.Item(49, 13).Value = -33264  # WVR!M49 (new cell)
$ws.Cells.Item(62, 8).Value = 16673666  # WVR!H62 was 18188864
$ws.Cells.Item(62, 10).Value = 18188818  # WVR!J62 was 20007050
$ws.Cells.Item(62, 12).Value = 18188818  # WVR!L62 was 20007050
$ws.Cells.Item(62, 14).Value = -18190066  # WVR!N62 was -20008298
$ws.Cells.Item(65, 8).Value = 16673666  # WVR!H65 was 18188864
$ws.Cells.Item(65, 10).Value = 18188818  # WVR!J65 was 20007050
$ws.Cells.Item(65, 12).Value = 90944090  # WVR!L65 was 100035250
$ws.Cells.Item(65, 14).Value = -90950330  # WVR!N65 was -100041490
$ws.Cells.Item(100, 8).Value = 973.6667  # WVR!H100 was 633.65216
$ws.Cells.Item(100, 9).Value = 1266.6  # WVR!I100 was 430.22223
$ws.Cells.Item(100, 11).Value = 2533.2  # WVR!K100 was 860.44446
$ws.Cells.Item(100, 13).Value = -1992.2  # WVR!M100 was -319.44446
